$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add artificial data row (row 8): a "만랩커피" expense entry with receipt checked
$ws.Range("B8").Value = "만랩커피"
$ws.Range("D8").Value = 11200
$ws.Range("E8").Value = "V"

# Update the active selection to reflect where the user ended up
$ws.Range("I12").Select()
